$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- ALC ---
# Row 6: Antidote
$ws1.Range("H6").Value = 15873245
$ws1.Range("I6").Value = 30303154
$ws1.Range("K6").Value = 90909462
$ws1.Range("M6").Value = -90909350

# Row 15: Ether
$ws1.Range("H15").Value = 1012.0714
$ws1.Range("I15").Value = 1012.0714
$ws1.Range("K15").Value = 3036.2142
$ws1.Range("M15").Value = -2867.2142

# Row 19: Roof Tile
$ws1.Range("H19").Value = 873.9
$ws1.Range("I19").Value = 422.33334
$ws1.Range("K19").Value = 422.33334
$ws1.Range("M19").Value = -247.33334

# Row 86: Enchanted Aurum Regis Ink
$ws1.Range("H86").Value = 50111
$ws1.Range("I86").Value = 4000
$ws1.Range("J86").Value = 86999.8
$ws1.Range("K86").Value = 4000
$ws1.Range("L86").Value = 86999.8
$ws1.Range("M86").Value = -2877
$ws1.Range("N86").Value = -89245.8

# Row 89: Enchanted Aurum Regis Ink
$ws1.Range("H89").Value = 50111
$ws1.Range("I89").Value = 4000
$ws1.Range("J89").Value = 86999.8
$ws1.Range("K89").Value = 20000
$ws1.Range("L89").Value = 434999
$ws1.Range("M89").Value = -14384
$ws1.Range("N89").Value = -446231

# Row 98: Enchanted Durium Ink
$ws1.Range("H98").Value = 61636.816
$ws1.Range("I98").Value = 62800.5
$ws1.Range("K98").Value = 62800.5
$ws1.Range("M98").Value = -61302.5

# Row 101: Cunning Craftsman's Tea
$ws1.Range("H101").Value = 9525835
$ws1.Range("I101").Value = 10206145
$ws1.Range("K101").Value = 30618435
$ws1.Range("M101").Value = -30616813

# Row 106: Enchanted Palladium Ink
$ws1.Range("H106").Value = 6294.0557
$ws1.Range("I106").Value = 8099.75
$ws1.Range("J106").Value = 4849.5
$ws1.Range("K106").Value = 8099.75
$ws1.Range("L106").Value = 4849.5
$ws1.Range("M106").Value = -7468.75
$ws1.Range("N106").Value = -6111.5

# Row 116: Growth Formula Kappa
$ws1.Range("H116").Value = 4705036.5
$ws1.Range("I116").Value = 4705036.5
$ws1.Range("J116").Value = 0
$ws1.Range("K116").Value = 4705036.5
$ws1.Range("L116").Value = 0
$ws1.Range("M116").Value = $null
$ws1.Range("N116").Value = -4701594.5

# Row 121: Tincture of Mind
$ws1.Range("H121").Value = 1726.5555
$ws1.Range("J121").Value = 1726.5555
$ws1.Range("L121").Value = 5179.666499999999
$ws1.Range("N121").Value = -8673.666499999999

# Row 122: Enchanted High Durium Ink
$ws1.Range("H122").Value = 61636.816
$ws1.Range("I122").Value = 62800.5
$ws1.Range("K122").Value = 188401.5
$ws1.Range("M122").Value = -185951.5

# Row 138: Cunning Craftsman's Tisane
$ws1.Range("H138").Value = 12624.95
$ws1.Range("I138").Value = 13714.286
$ws1.Range("J138").Value = 12038.385
$ws1.Range("K138").Value = 41142.858
$ws1.Range("L138").Value = 36115.155
$ws1.Range("M138").Value = -36002.858
$ws1.Range("N138").Value = -46395.155

# Row 141: Grade 1 Gemdraught of Mind
$ws1.Range("H141").Value = 6971.0527
$ws1.Range("I141").Value = 4769.909
$ws1.Range("K141").Value = 14309.727
$ws1.Range("M141").Value = -9129.726999999999

# --- ARM ---
# Row 32: Steel Ingot
$ws2.Range("H32").Value = 2719.45
$ws2.Range("I32").Value = 2811.8108
$ws2.Range("J32").Value = 1580.3334
$ws2.Range("K32").Value = 2811.8108
$ws2.Range("L32").Value = 1580.3334
$ws2.Range("M32").Value = -2524.8108
$ws2.Range("N32").Value = -2154.3334

# Row 43: Steel Sabatons
$ws2.Range("H43").Value = 11890.3
$ws2.Range("I43").Value = 8315.333000000001
$ws2.Range("J43").Value = 13422.429
$ws2.Range("K43").Value = 8315.333000000001
$ws2.Range("L43").Value = 13422.429
$ws2.Range("M43").Value = -8002.333000000001
$ws2.Range("N43").Value = -14048.429

# Row 44: Mythril Plate
$ws2.Range("H44").Value = 80706.14
$ws2.Range("J44").Value = 80706.14
$ws2.Range("L44").Value = 80706.14
$ws2.Range("N44").Value = -81682.14

# Row 45: Mythril Ingot
$ws2.Range("H45").Value = 7440.1724
$ws2.Range("I45").Value = 12016.267
$ws2.Range("J45").Value = 2537.2144
$ws2.Range("K45").Value = 12016.267
$ws2.Range("L45").Value = 2537.2144
$ws2.Range("M45").Value = -11639.267
$ws2.Range("N45").Value = -3291.2144

# Row 46: Heavy Steel Flanchard
$ws2.Range("H46").Value = 11028.8
$ws2.Range("I46").Value = 0
$ws2.Range("K46").Value = 0
$ws2.Range("M46").Value = $null

# Row 61: Cobalt Ingot
$ws2.Range("H61").Value = 2485.4666
$ws2.Range("I61").Value = 928.5
$ws2.Range("J61").Value = 5599.4
$ws2.Range("K61").Value = 928.5
$ws2.Range("L61").Value = 5599.4
$ws2.Range("M61").Value = -716.5
$ws2.Range("N61").Value = -6023.4

# Row 74: Titanium Nugget
$ws2.Range("H74").Value = 77280.84
$ws2.Range("I74").Value = 85813.88
$ws2.Range("J74").Value = 4750
$ws2.Range("K74").Value = 85813.88
$ws2.Range("L74").Value = 4750
$ws2.Range("M74").Value = -84939.88
$ws2.Range("N74").Value = -6498

# Row 77: Titanium Nugget
$ws2.Range("H77").Value = 77280.84
$ws2.Range("I77").Value = 85813.88
$ws2.Range("J77").Value = 4750
$ws2.Range("K77").Value = 429069.4
$ws2.Range("L77").Value = 23750
$ws2.Range("M77").Value = -424701.4
$ws2.Range("N77").Value = -32486

# Row 97: High Steel Ingot
$ws2.Range("H97").Value = 4678.52
$ws2.Range("I97").Value = 5207.7617
$ws2.Range("K97").Value = 5207.7617
$ws2.Range("M97").Value = -4711.7617

# Row 132: Mountain Chromite Ingot
$ws2.Range("H132").Value = 2691.3076
$ws2.Range("I132").Value = 2153.2144
$ws2.Range("J132").Value = 4951.3
$ws2.Range("K132").Value = 6459.6432
$ws2.Range("L132").Value = 14853.9
$ws2.Range("M132").Value = -3929.6432
$ws2.Range("N132").Value = -19913.9

# Row 136: Cobalt Tungsten Ingot
$ws2.Range("H136").Value = 2485.4666
$ws2.Range("I136").Value = 928.5
$ws2.Range("J136").Value = 5599.4
$ws2.Range("K136").Value = 2785.5
$ws2.Range("L136").Value = 16798.2
$ws2.Range("M136").Value = -235.5
$ws2.Range("N136").Value = -21898.2

# --- BSM ---
# Row 102: Doman Steel Mortar
$ws3.Range("H102").Value = 16499.8
$ws3.Range("I102").Value = 11874.75
$ws3.Range("J102").Value = 35000
$ws3.Range("K102").Value = 11874.75
$ws3.Range("L102").Value = 35000
$ws3.Range("M102").Value = -8629.75
$ws3.Range("N102").Value = -41490

# Row 107: Deepgold Nugget
$ws3.Range("H107").Value = 2022.579
$ws3.Range("I107").Value = 1875.6471
$ws3.Range("K107").Value = 1875.6471
$ws3.Range("M107").Value = 44.35290000000009

# Row 134: Ruthenium Ingot
$ws3.Range("H134").Value = 7333.5713
$ws3.Range("I134").Value = 7744.8335
$ws3.Range("J134").Value = 4866
$ws3.Range("K134").Value = 23234.5005
$ws3.Range("L134").Value = 14598
$ws3.Range("M134").Value = -20699.5005
$ws3.Range("N134").Value = -19668

# --- CRP ---
# Row 18: Ash Spinning Wheel
$ws4.Range("H18").Value = 71012.5
$ws4.Range("J18").Value = 71012.5
$ws4.Range("L18").Value = 71012.5
$ws4.Range("N18").Value = -71472.5

# Row 31: Walnut Lumber
$ws4.Range("H31").Value = 1917.3732
$ws4.Range("I31").Value = 1216.4445
$ws4.Range("J31").Value = 2731.3547
$ws4.Range("K31").Value = 1216.4445
$ws4.Range("L31").Value = 2731.3547
$ws4.Range("M31").Value = -921.4445000000001
$ws4.Range("N31").Value = -3321.3547

# Row 34: Walnut Lumber
$ws4.Range("H34").Value = 1917.3732
$ws4.Range("I34").Value = 1216.4445
$ws4.Range("J34").Value = 2731.3547
$ws4.Range("K34").Value = 1216.4445
$ws4.Range("L34").Value = 2731.3547
$ws4.Range("M34").Value = -1014.4445
$ws4.Range("N34").Value = -3135.3547

# Row 58: Mahogany Lumber
$ws4.Range("H58").Value = 6764.467
$ws4.Range("I58").Value = 15034.2
$ws4.Range("J58").Value = 2629.6
$ws4.Range("K58").Value = 15034.2
$ws4.Range("L58").Value = 2629.6
$ws4.Range("M58").Value = -14831.2
$ws4.Range("N58").Value = -3035.6

# Row 59: Crab Bow
$ws4.Range("H59").Value = 56343.625
$ws4.Range("J59").Value = 58107
$ws4.Range("L59").Value = 58107
$ws4.Range("N59").Value = -60397

# Row 74: Dark Chestnut Rod
$ws4.Range("H74").Value = 41786.855
$ws4.Range("J74").Value = 44744.8
$ws4.Range("L74").Value = 44744.8
$ws4.Range("N74").Value = -46492.8

# Row 77: Dark Chestnut Rod
$ws4.Range("H77").Value = 41786.855
$ws4.Range("J77").Value = 44744.8
$ws4.Range("L77").Value = 134234.4
$ws4.Range("N77").Value = -142970.4

# Row 94: Beech Lumber
$ws4.Range("H94").Value = 2446.9333
$ws4.Range("I94").Value = 5670.5
$ws4.Range("J94").Value = 1274.7273
$ws4.Range("K94").Value = 5670.5
$ws4.Range("L94").Value = 1274.7273
$ws4.Range("M94").Value = -5219.5
$ws4.Range("N94").Value = -2176.7273

# Row 136: Dark Mahogany Lumber
$ws4.Range("H136").Value = 6764.467
$ws4.Range("I136").Value = 15034.2
$ws4.Range("J136").Value = 2629.6
$ws4.Range("K136").Value = 45102.60000000001
$ws4.Range("L136").Value = 7888.799999999999
$ws4.Range("M136").Value = -42552.60000000001
$ws4.Range("N136").Value = -12988.8

# --- CUL ---
# Row 139: Wild Banana Blend
$ws5.Range("H139").Value = 1877796.6
$ws5.Range("I139").Value = 2729040.8
$ws5.Range("J139").Value = 5059.8
$ws5.Range("K139").Value = 8187122.399999999
$ws5.Range("L139").Value = 15179.4
$ws5.Range("M139").Value = -8181982.399999999
$ws5.Range("N139").Value = -25459.4

# --- GSM ---
# Row 132: Lar Ingot
$ws6.Range("H132").Value = 4019.2
$ws6.Range("I132").Value = 3525.261
$ws6.Range("J132").Value = 9699.5
$ws6.Range("K132").Value = 10575.783
$ws6.Range("L132").Value = 29098.5
$ws6.Range("M132").Value = -8045.782999999999
$ws6.Range("N132").Value = -34158.5

# --- LTW ---
# Row 22: Aldgoat Leather
$ws7.Range("H22").Value = 6986.5938
$ws7.Range("I22").Value = 8595.5
$ws7.Range("J22").Value = 4918
$ws7.Range("K22").Value = 8595.5
$ws7.Range("L22").Value = 4918
$ws7.Range("M22").Value = -8300.5
$ws7.Range("N22").Value = -5508

# Row 27: Aldgoat Leather
$ws7.Range("H27").Value = 6986.5938
$ws7.Range("I27").Value = 8595.5
$ws7.Range("J27").Value = 4918
$ws7.Range("K27").Value = 8595.5
$ws7.Range("L27").Value = 4918
$ws7.Range("M27").Value = -8488.5
$ws7.Range("N27").Value = -5132

# Row 141: Gargantuaskin Trousers of Striking
$ws7.Range("H141").Value = 49000
$ws7.Range("J141").Value = 49000
$ws7.Range("L141").Value = 49000
$ws7.Range("N141").Value = -59360

# --- WVR ---
# Row 5: Hempen Halfgloves
$ws8.Range("H5").Value = 22012500
$ws8.Range("I5").Value = 18014286
$ws8.Range("K5").Value = 18014286
$ws8.Range("M5").Value = -18014174

# Row 96: Ruby Cotton Cloth
$ws8.Range("H96").Value = 6667830.5
$ws8.Range("I96").Value = 9091979
$ws8.Range("J96").Value = 1423.5
$ws8.Range("K96").Value = 9091979
$ws8.Range("L96").Value = 1423.5
$ws8.Range("M96").Value = -9090606
$ws8.Range("N96").Value = -4169.5

# Row 132: Snow Cotton Cloth
$ws8.Range("H132").Value = 15531.857
$ws8.Range("I132").Value = 18201.969
$ws8.Range("J132").Value = 5741.4443
$ws8.Range("K132").Value = 54605.90700000001
$ws8.Range("L132").Value = 17224.3329
$ws8.Range("M132").Value = -52075.90700000001
$ws8.Range("N132").Value = -22284.3329
